# Apply the "assertMatch(text,regex)" base-command addition (and the
# bundled openFile(filePath)/tn.5250-removal housekeeping) to the hidden
# '#system' reference sheet, mirroring the target OOXML diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# ---------------------------------------------------------------------
# 1) Column F ("base" named range): insert "assertMatch(text,regex)" at
#    F11 (alphabetically between assertEqual and assertNotContain),
#    pushing F11:F44 down to F12:F45.
# ---------------------------------------------------------------------
for ($r = 44; $r -ge 11; $r--) {
    $v = $ws.Cells.Item($r, 6).Text
    $ws.Cells.Item($r + 1, 6).Value = $v
}
$ws.Cells.Item(11, 6).Value = "assertMatch(text,regex)"

# ---------------------------------------------------------------------
# 2) Column J ("external" named range): insert "openFile(filePath)" at
#    J2 (alphabetically before runJUnit), pushing J2:J6 down to J3:J7.
# ---------------------------------------------------------------------
for ($r = 6; $r -ge 2; $r--) {
    $v = $ws.Cells.Item($r, 10).Text
    $ws.Cells.Item($r + 1, 10).Value = $v
}
$ws.Cells.Item(2, 10).Value = "openFile(filePath)"

# ---------------------------------------------------------------------
# 3) Column A ("target" named range): remove the "tn.5250" entry (row
#    27), pulling A28:A33 up to A27:A32.
# ---------------------------------------------------------------------
for ($r = 27; $r -le 32; $r++) {
    $v = $ws.Cells.Item($r + 1, 1).Text
    $ws.Cells.Item($r, 1).Value = $v
}
$ws.Cells.Item(33, 1).Value = ""

# ---------------------------------------------------------------------
# 4) Columns AA:AG: the "tn.5250" command list (old column AA) is
#    removed, so "web"/"webalert"/"webcookie"/"ws"/"ws.async"/"xml"
#    each shift one column to the left (AB->AA, AC->AB, AD->AC,
#    AE->AD, AF->AE, AG->AF), for every row that holds data (1-151).
# ---------------------------------------------------------------------
for ($r = 1; $r -le 151; $r++) {
    for ($c = 27; $c -le 32; $c++) {
        $v = $ws.Cells.Item($r, $c + 1).Text
        $ws.Cells.Item($r, $c).Value = $v
    }
    $ws.Cells.Item($r, 33).Value = ""
}

# ---------------------------------------------------------------------
# 5) Fix up the defined names that reference the shifted ranges.
#    ("tn.5250" itself is intentionally left untouched/orphaned, as in
#    the target diff.)
# ---------------------------------------------------------------------
$wb.Names.Item("base").RefersTo        = "='#system'!`$F`$2:`$F`$45"
$wb.Names.Item("external").RefersTo    = "='#system'!`$J`$2:`$J`$7"
$wb.Names.Item("target").RefersTo      = "='#system'!`$A`$2:`$A`$32"
$wb.Names.Item("web").RefersTo         = "='#system'!`$AA`$2:`$AA`$151"
$wb.Names.Item("webalert").RefersTo    = "='#system'!`$AB`$2:`$AB`$8"
$wb.Names.Item("webcookie").RefersTo   = "='#system'!`$AC`$2:`$AC`$10"
$wb.Names.Item("ws").RefersTo          = "='#system'!`$AD`$2:`$AD`$17"
$wb.Names.Item("ws.async").RefersTo    = "='#system'!`$AE`$2:`$AE`$8"
$wb.Names.Item("xml").RefersTo         = "='#system'!`$AF`$2:`$AF`$27"
